# Repull data, push all data, mean calculation
# Updates the "dSF" column (column F) values for a set of rows to reflect
# re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -2
    5  = 2
    7  = 2
    8  = 0
    10 = -2
    17 = 3
    18 = -1
    25 = -4
    26 = -3
    27 = 2
    29 = -2
    30 = -2
    31 = -2
    32 = -6
    48 = -4
    49 = -1
    51 = 2
    53 = -6
    57 = -3
    58 = 4
    61 = 2
    63 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
